# Update column G ("K") values for rows 2-8 per regenerated save_data
# (K column now derived from strikeout-based stat instead of Strike# count)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 1
$ws.Range("G8").Value = 3
